$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching the style of the other
# header cells (e.g. G1: bold, centered, bordered)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add new "Save" column data values
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
